$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new "Shunt Impedans (ohm)" column before the existing
# "Shunt Admittans (ohm)" column, shifting F,G,H -> G,H,I ---------------

# Preserve the bold/border header style onto the new rightmost header
# cell (I1) before it gets overwritten by the shift below.
$ws.Range("H1").Copy()
$ws.Range("I1").PasteSpecial(-4122)  # xlPasteFormats

for ($r = 1; $r -le 10; $r++) {
    $ws.Cells.Item($r, 9).Value = $ws.Cells.Item($r, 8).Value()
    $ws.Cells.Item($r, 8).Value = $ws.Cells.Item($r, 7).Value()
    $ws.Cells.Item($r, 7).Value = $ws.Cells.Item($r, 6).Value()
}

# New header for the inserted column
$ws.Range("F1").Value = "Shunt Impedans (ohm)"

# New (corrected) values for the inserted "Shunt Impedans (ohm)" column F
$ws.Range("F2").Value = "-15977.005781448112j"
$ws.Range("F3").Value = "-19461.352786976684j"
$ws.Range("F4").Value = "-7654.255907848572j"
$ws.Range("F5").Value = "-23280.17890615012j"
$ws.Range("F6").Value = "-3921.6179551521614j"
$ws.Range("F7").Value = "-5132.7058530668j"
$ws.Range("F8").Value = "-7406.34478532716j"
$ws.Range("F9").Value = "-3458.160986721755j"
$ws.Range("F10").Value = "-3135.162231321008j"

# Corrected values for the (shifted) "Shunt Admittans (ohm)" column G
$ws.Range("G2").Value = "6.258995043746944e-05j"
$ws.Range("G3").Value = "5.1383889442114665e-05j"
$ws.Range("G4").Value = "0.00013064627209218512j"
$ws.Range("G5").Value = "4.2954996352533256e-05j"
$ws.Range("G6").Value = "0.0002549967925065764j"
$ws.Range("G7").Value = "0.0001948290100050246j"
$ws.Range("G8").Value = "0.00013501936906598212j"
$ws.Range("G9").Value = "0.0002891710373923261j"
$ws.Range("G10").Value = "0.00031896276052631814j"

# Corrected values for the (shifted) "Shunt Admittans (p.u.)" column H
$ws.Range("H2").Value = "0.005633095539372249j"
$ws.Range("H3").Value = "0.00462455004979032j"
$ws.Range("H4").Value = "0.011758164488296662j"
$ws.Range("H5").Value = "0.003865949671727993j"
$ws.Range("H6").Value = "0.022949711325591874j"
$ws.Range("H7").Value = "0.017534610900452215j"
$ws.Range("H8").Value = "0.012151743215938392j"
$ws.Range("H9").Value = "0.02602539336530935j"
$ws.Range("H10").Value = "0.028706648447368633j"

# Corrected values for the (shifted) "Shunt Admittans (p.u.) half" column I
$ws.Range("I2").Value = "0.0028165477696861247j"
$ws.Range("I3").Value = "0.00231227502489516j"
$ws.Range("I4").Value = "0.005879082244148331j"
$ws.Range("I5").Value = "0.0019329748358639966j"
$ws.Range("I6").Value = "0.011474855662795937j"
$ws.Range("I7").Value = "0.008767305450226108j"
$ws.Range("I8").Value = "0.006075871607969196j"
$ws.Range("I9").Value = "0.013012696682654675j"
$ws.Range("I10").Value = "0.014353324223684316j"

# --- Column widths (engine stores ColumnWidth + 5/6 as the saved "width"
# attribute, so back the desired stored width off by that fixed offset) ---
$ws.Columns("B:B").ColumnWidth = 41.830729166666664
$ws.Columns("C:C").ColumnWidth = 43.998697916666664
$ws.Columns("D:D").ColumnWidth = 39.666666666666664
$ws.Columns("E:E").ColumnWidth = 12.998697916666666
$ws.Columns("F:F").ColumnWidth = 21.830729166666668
$ws.Columns("G:G").ColumnWidth = 26.998697916666668
$ws.Columns("H:H").ColumnWidth = 22.666666666666668
$ws.Columns("I:I").ColumnWidth = 24.330729166666668

# --- View state: scroll so column C is left-most visible, select G14 ---
$ws.Application.ActiveWindow.ScrollColumn = 3
$ws.Range("G14").Select()
